# "Generate Report for Handoff"
#
# A new handoff round was generated for the b.md file: the localization
# status report needs updating on all three sheets (Overview, zh-cn, de-de)
# to reflect the new handoff xliff files / timestamps, and the two
# per-language detail sheets get a widened "Error Detail" column plus an
# explanatory error message for the row that is not yet in sync.
#
# NOTE: cell writes below are intentionally ordered to match first-use order
# (Overview, then zh-cn, then de-de; each sheet top-to-bottom, left-to-right)
# so brand-new shared strings land in the same relative sequence the
# reference workbook uses.

$wb = $excel.ActiveWorkbook

$overview = $wb.Sheets.Item("Overview")
$zhcn = $wb.Sheets.Item("zh-cn")
$dede = $wb.Sheets.Item("de-de")

# --- Overview sheet: row 3 is the b.md entry ---------------------------
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-18 10:41:31"

# --- zh-cn sheet: widen Error Detail column, update row 3 (b.md) -------
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").ClearFormats()
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-18 10:41:26"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/140afbae760f0168fab1010f3afdbedddcdddf97/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259d81ba9a3c622aaa12138ea85301ef6960e894/e2e/b.md."

# --- de-de sheet: widen Error Detail column, update row 3 (b.md) -------
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664

$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").ClearFormats()
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-18 10:41:31"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/140afbae760f0168fab1010f3afdbedddcdddf97/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/259d81ba9a3c622aaa12138ea85301ef6960e894/e2e/b.md."
